$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aufgaben")

# 1. Update the "Letzte Aktualisierung" date in A5
$ws.Range("A5").Value = "Letzte Aktualisierung: 10.07.2015"

# 2. Update row 24's task name (Piwik-Analyse), rename to include "Seitenansichten:"
$ws.Range("A24").Value = "Piwik-Analyse (Seitenansichten: Ebene 1 +2)"

# 3. Insert two new rows after row 24 (new rows 25 and 26), pushing old rows 25-27 down to 27-29
$ws.Rows.Item(25).Insert()
$ws.Rows.Item(25).Insert()

# Fill new row 25: Analyse des SRS-Templates
$ws.Range("A25").Value = "Analyse des SRS-Templates"
$ws.Range("B25").Value = "Benedikt Häring, Johannes Vogl "
$ws.Range("C25").Value = 1/1440
$ws.Range("D25").Value = "SRS-Dokument erstellen"

# Fill new row 26: Besprechung der naechsten Projektschritte
$ws.Range("A26").Value = "Besprechung der nächsten Projektschritte "
$ws.Range("B26").Value = "Alle Teammitglieder beteiligt"
$ws.Range("C26").Value = 1/1440
$ws.Range("D26").Value = "Projektplanung"

# 4. Insert four new rows after row 29 (the old row27 "Contextual Inquiry Videoauswertung" now at 29)
$ws.Rows.Item(30).Insert()
$ws.Rows.Item(30).Insert()
$ws.Rows.Item(30).Insert()
$ws.Rows.Item(30).Insert()

# Fill new row 30: Meeting mit Raphael Wimmer
$ws.Range("A30").Value = "Meeting mit Raphael Wimmer"
$ws.Range("B30").Value = "Alle Teammitglieder beteiligt"
$ws.Range("C30").Value = 0.5/1440
$ws.Range("D30").Value = "Projektplanung"

# Fill new row 31: Content Analyse der RZ-Website Teil 2 (no duration value)
$ws.Range("A31").Value = "Content Analyse der RZ-Website Teil 2"
$ws.Range("B31").Value = "Fabian Huth"
$ws.Range("D31").Value = "Analyse des Ist-Zustands"
$ws.Range("C31").Clear()

# Fill new row 32: Piwik-Analyse (Seitenansichten: Ebene 3)
$ws.Range("A32").Value = "Piwik-Analyse (Seitenansichten: Ebene 3)"
$ws.Range("B32").Value = "Dominik Bauer"
$ws.Range("C32").Value = 1.5/1440
$ws.Range("D32").Value = "Analyse des Ist-Zustands"

# Fill new row 33: Piwik-Analyse (Seitenansichten Visualisierung)
$ws.Range("A33").Value = "Piwik-Analyse (Seitenansichten Visualisierung) "
$ws.Range("B33").Value = "Dominik Bauer"
$ws.Range("C33").Value = 1/1440
$ws.Range("D33").Value = "Analyse des Ist-Zustands"

# 5. Update sheet view: zoom 70 -> 60, pane topLeftCell A15 -> A22, selection A5 -> D30
$ws.Activate()
$excel.ActiveWindow.Zoom = 60
$ws.Range("D30").Select()
$excel.ActiveWindow.ScrollRow = 22
